$d = $word.ActiveDocument

# --- Change 1: merge the two runs of the "22/09- ..." paragraph into a single run ---
$null = $d.Content.Find.Execute("22/09- Definição de quem ficou responsável pela atualização das atas de reunião, trello e quadro.", $true, $false, $false, $false, $false, $true, 1, $false, "22/09- Definição de quem ficou responsável pela atualização das atas de reunião, trello e quadro.", 2)

# --- Change 2: locate the "22/09-..." paragraph and the empty paragraph right after it ---
$targetPara = $null
$nextPara = $null
foreach ($p in $d.Paragraphs) {
    if ($targetPara -ne $null) {
        $nextPara = $p
        break
    }
    if ($p.Range.Text -like "22/09-*") {
        $targetPara = $p
    }
}

# Fill the next (currently empty) paragraph with the "25/09" minute entry, reusing its run
$newText = "25/09 – Definição das atividades da semana do dia 27/09 á 02/10. "
$r = $nextPara.Range
$r.Text = $newText

# Re-fetch paragraph (ranges can re-seat after edits) and set explicit run formatting
$nextPara2 = $targetPara.Next()
$textRange = $d.Range($nextPara2.Range.Start, $nextPara2.Range.Start + $newText.Length)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Arial" w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr><w:t xml:space="preserve">' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $textRange.InsertXML($xml)

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
